$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look like pure numbers are stored as text,
# matching the source data which treats these as display strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '72.166.74'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '4.037.29'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '531.77'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '151.62'
$ws.Range('E6').Value = '  +2.00%  '
$ws.Range('D7').Value = '0.698'
$ws.Range('E7').Value = '  +11.71%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.751'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D11').Value = '0.0000328'
$ws.Range('E11').Value = '  -3.28%  '
$ws.Range('D12').Value = '47.97'
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.692.58'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '10.68'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '4.044.72'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').Value = '14.15'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '20.61'
$ws.Range('E17').Value = '  -3.82%  '
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D20').Value = '72.155.74'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D21').Value = '429.07'
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('D22').Value = '98.14'
$ws.Range('E22').Value = '  +4.02%  '
$ws.Range('D23').Value = '3.49'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('E24').Value = '  +4.40%  '
$ws.Range('D25').Value = '14.34'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '11.18'
$ws.Range('E26').Value = '  -9.62%  '
$ws.Range('D27').Value = '10.76'
$ws.Range('E27').Value = '  -4.17%  '
$ws.Range('D28').Value = '5.84'
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('D29').Value = '36.83'
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('D30').Value = '3.58'
$ws.Range('E30').Value = '  +21.97%  '
$ws.Range('D31').Value = '13.42'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  -1.85%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '7.14'
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').Value = '676.53'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('D35').Value = '44.51'
$ws.Range('E35').Value = '  +9.02%  '
$ws.Range('D36').Value = '66.15'
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('D37').Value = '0.447'
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('D38').Value = '0.0₃0834'
$ws.Range('E38').Value = '  -8.14%  '
$ws.Range('E39').Value = '  -2.88%  '
$ws.Range('D40').Value = '3.39'
$ws.Range('E40').Value = '  -4.92%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '0.0487'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('D44').Value = '3.19'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '9.71'
$ws.Range('E46').Value = '  +5.87%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.43'
$ws.Range('E47').Value = '  -3.61%  '
$ws.Range('D48').Value = '2.63'
$ws.Range('E48').Value = '  -6.32%  '
$ws.Range('D49').Value = '3.02'
$ws.Range('E49').Value = '  -6.01%  '
$ws.Range('E50').Value = '  -3.26%  '
$ws.Range('D51').Value = '145.73'
$ws.Range('E51').Value = '  +1.25%  '
